# weight_tracker: "add updates and weekly trend tracker"
# Append two new weigh-ins (rows 71 & 72) to the raw_data sheet, matching
# the existing table's layout: date+time (col A), time-of-day fraction
# (col B), weight (col C), and a shared AM/PM formula (col D).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("raw_data")

# Number formats used by the existing rows (style indices 6 and 2 in the
# original workbook: numFmtId 22 "m/d/yy h:mm" for the date+time column,
# numFmtId 20 "h:mm" for the time-of-day column).
$dateTimeFormat = "m/d/yy h:mm"
$timeFormat = "h:mm"

# Row 71: 2020-08-31 21:05 -> 74.7 lbs, PM
$ws.Range("A71").Value = 44074.878472222219
$ws.Range("A71").NumberFormat = $dateTimeFormat
$ws.Range("B71").Value = 0.87847222222222221
$ws.Range("B71").NumberFormat = $timeFormat
$ws.Range("C71").Value = 74.7
$ws.Range("D71").Formula = "=IF(B71<TIME(12,0,0), ""AM"", ""PM"")"

# Row 72: 2020-09-01 08:16 -> 73 lbs, AM
$ws.Range("A72").Value = 44075.344444444447
$ws.Range("A72").NumberFormat = $dateTimeFormat
$ws.Range("B72").Value = 0.3444444444444445
$ws.Range("B72").NumberFormat = $timeFormat
$ws.Range("C72").Value = 73
$ws.Range("D72").Formula = "=IF(B72<TIME(12,0,0), ""AM"", ""PM"")"

# Match the author's final selection/cursor position on the new last row.
$ws.Range("A72").Select() | Out-Null
